$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("889104", "Guy", "Dumas"),
    @("433402", "Ayn", "Sawyer"),
    @("216873", "Guy", "Holmes"),
    @("981322", "Isaac", "Montag"),
    @("562190", "Tom", "Holmes"),
    @("593299", "Sherlock", "Alighieri"),
    @("830210", "Tom", "Montag"),
    @("539227", "Guy", "Sawyer"),
    @("350556", "William", "Sanderson"),
    @("244761", "Alexandre", "Sawyer")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    # Column A holds numeric-looking IDs that must stay text, like the
    # source workbook (openpyxl inline strings). Force text format first
    # so Excel doesn't auto-coerce them into numbers.
    $cellA = $ws.Cells.Item($row, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
